{"js": "// Wrap the address-format placeholder text in curly quotes (\u201c \u2026 \u201d) in the\n// four table cells that describe the required address format, e.g.\n//   <denominazione urbanistica>, <indirizzo>, <numero civico>, <CAP>, <Comune>, <Provincia>.\n// becomes\n//   \"<denominazione urbanistica>, <indirizzo>, <numero civico>, <CAP>, <Comune>, <Provincia>\".\n// (using the Italian/typographic double quotes U+201C / U+201D).\n\nconst oldText =\n  \"<denominazione urbanistica>, <indirizzo>, <numero civico>, <CAP>, <Comune>, <Provincia>.\";\nconst newText =\n  \"\\u201C<denominazione urbanistica>, <indirizzo>, <numero civico>, <CAP>, <Comune>, <Provincia>\\u201D.\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (const range of results.items) {\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Wrap the address-format placeholder text in curly quotes (\"...\" using the\n# typographic double quotes U+201C / U+201D) in the four table cells that\n# describe the required address format, e.g.\n#   <denominazione urbanistica>, <indirizzo>, <numero civico>, <CAP>, <Comune>, <Provincia>.\n# becomes\n#   \"<denominazione urbanistica>, <indirizzo>, <numero civico>, <CAP>, <Comune>, <Provincia>\".\n\n$d = $word.ActiveDocument\n\n$openQuote  = [char]0x201C\n$closeQuote = [char]0x201D\n\n$oldText = \"<denominazione urbanistica>, <indirizzo>, <numero civico>, <CAP>, <Comune>, <Provincia>.\"\n$newText = $openQuote + \"<denominazione urbanistica>, <indirizzo>, <numero civico>, <CAP>, <Comune>, <Provincia>\" + $closeQuote + \".\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n\n# MatchCase:=True, MatchWholeWord:=False, MatchWildcards:=False,\n# MatchSoundsLike:=False, MatchAllWordForms:=False, Forward:=True,\n# Wrap:=wdFindContinue(1), Format:=False, Replace:=wdReplaceAll(2)\n$find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n"}
